$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room: insert 2 rows above the old row 13 so that the
#    existing "stripe" / "vercel" credential blocks (old rows 13-17)
#    shift down to rows 15-19, matching the target layout.
# ------------------------------------------------------------------
$ws.Rows("13:14").Insert()

# ------------------------------------------------------------------
# 2) New row 13: "environment variables" note block.
# ------------------------------------------------------------------
$ws.Range("B13").Value = "環境変数の設定"
$ws.Range("C13").Value = "settingから行うが、チーム画面に勝手になるので超注意。`nSettingボタンを押し、左上の上と下の記号をおして、task-managerを選択するとProjectSettingの表示になる。Settingだとだめ"
$ws.Range("C13").WrapText = $true
$ws.Rows("13:13").RowHeight = 58.5

# ------------------------------------------------------------------
# 3) New test-account blocks at rows 21-31 (3 blocks of 3 rows,
#    separated by blank rows, mirroring the stripe/vercel layout).
# ------------------------------------------------------------------

# Block 1 - rows 21-23
$ws.Range("A21").Value = "テストID"
$ws.Range("B21").Value = "めーる"
$ws.Range("C21").Value = "zxjfurhjvc@yahoo.co.jp"
$ws.Range("B22").Value = "ID"
$ws.Range("C22").Value = "bbbb"
$ws.Range("B23").Value = "Pass"
$ws.Range("C23").Value = "dgfhyt6543g"

# Block 2 - rows 25-27
$ws.Range("B25").Value = "めーる"
$ws.Range("C25").Value = "vndhry@yahoo.co.jp"
$ws.Range("B26").Value = "ID"
$ws.Range("C26").Value = "CCCC"
$ws.Range("B27").Value = "Pass"
$ws.Range("C27").Value = "dhfjt6534"

# Block 3 - rows 29-31
$ws.Range("B29").Value = "めーる"
$ws.Range("C29").Value = "shinji19750918@yahoo.co.jp"
$ws.Range("B30").Value = "ID"
$ws.Range("C30").Value = "shishi0918"
$ws.Range("B31").Value = "Pass"
$ws.Range("C31").Value = "shishi0918"

# ------------------------------------------------------------------
# 4) Apply cell formatting to match the other credential blocks by
#    copying the formats from the already-existing, equivalently
#    styled cells (format-only paste, values are untouched).
# ------------------------------------------------------------------
$xlPasteFormats = -4122

# style "s=3" cells (plain label style used throughout the ID/PASS rows)
$ws.Range("B16").Copy()
$ws.Range("B21,B22,B23,B25,B26,B27,B29,B30,B31").PasteSpecial($xlPasteFormats)

# style "s=4" cells (hyperlink style variant used for stripe/url hyperlink cells)
$ws.Range("C15").Copy()
$ws.Range("C21,C25").PasteSpecial($xlPasteFormats)

# style "s=6" cells (hyperlink style variant used for mailto hyperlink cells)
$ws.Range("C18").Copy()
$ws.Range("C29").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 5) Rebuild hyperlinks. Row-insert does not renumber the existing
#    hyperlink anchors, so clear them all and re-add at the correct,
#    final cell locations.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:shinji19750918@yahoo.co.jp")
$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:shinji19750918@yahoo.co.jp")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://supabase.com/dashboard/project/mwewfabykeyjopbeqbfj")
$ws.Hyperlinks.Add($ws.Range("C21"), "mailto:zxjfurhjvc@yahoo.co.jp")
$ws.Hyperlinks.Add($ws.Range("C25"), "mailto:vndhry@yahoo.co.jp")
$ws.Hyperlinks.Add($ws.Range("C29"), "mailto:shinji19750918@yahoo.co.jp")

# ------------------------------------------------------------------
# 6) Update the view selection to match.
# ------------------------------------------------------------------
$ws.Range("D30").Select()

Write-Host "done"
